$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row (row 7) mirroring the header row values: a new USERNAME/PASSWORD entry
$ws.Range("A7").Value = "USERNAME"
$ws.Range("B7").Value = "PASSWORD"
$ws.Range("C7").Value = 1

# Match the style used by the other data rows (style index 1 - Calibri 12, no special formatting)
$ws.Range("A7:C7").Font.Bold = $false

# Update the active selection to the newly added row, as in the target workbook
$ws.Range("F7").Select()
